# Fruta / hortaliza, semanal
#
# A new weekly record is inserted as row 123 (pushing the existing rows
# 123-287 down to 124-288, which is why the sheet's used range grows from
# A1:R287 to A1:R288). Populate the new row with the same constant
# attributes used by every other record in this sub-sheet, plus the new
# date (D) and volume (J) observations; price columns (K/L/M/N/P) keep
# their usual default values for a non-special Volumen reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 123..287 down one row, creating a blank row 123.
$ws.Rows.Item(123).Insert()

$ws.Range('A123').Value = 3
$ws.Range('B123').Value = 'Femacal de La Calera'
$ws.Range('C123').Value = 'Coquimbo'
$ws.Range('D123').Value = 44638
$ws.Range('E123').Value = 5
$ws.Range('F123').Value = 100112039
$ws.Range('G123').Value = 'Ciboulette'
$ws.Range('H123').Value = 'Sin especificar'
$ws.Range('I123').Value = 'Primera'
$ws.Range('J123').Value = 180
$ws.Range('K123').Value = 1500
$ws.Range('L123').Value = 1500
$ws.Range('M123').Value = 1500
$ws.Range('N123').Value = '$/docena de atados'
$ws.Range('O123').Value = 'Provincia de Quillota'
$ws.Range('P123').Value = 500
$ws.Range('Q123').Value = 3
$ws.Range('R123').Value = 'Hortaliza'
